$d = $word.ActiveDocument

# Locate the "Caching Implementation (Amazon S3)" bullet under "Project Goals:"
# (numId 4, ilvl 1) and the preceding sibling bullet so we can clone its
# paragraph formatting (style + list numbering) for the new bullet.
$target = $null
$donor = $null
$prev = $null
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "Caching Implementation (Amazon S3)") {
        $target = $p
        $donor = $prev
        break
    }
    $prev = $p
}

if ($target -ne $null) {
    # Insert a new empty paragraph right after the target bullet.
    $target.Range.InsertParagraphAfter()

    # The freshly inserted (still empty) paragraph now immediately follows
    # the target bullet; re-fetch it by index so we have live Start/End.
    $newIndex = $target.Index + 1
    $newPara = $d.Paragraphs($newIndex)

    # Clone the list/paragraph formatting from a sibling bullet at the same
    # list level (e.g. "Search for users and locations") ...
    if ($donor -ne $null) {
        $newPara.Range.FormattedText = $donor.Range.FormattedText
    }

    # ... then overwrite its text with the new bullet content.
    $newPara = $d.Paragraphs($newIndex)
    $newRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
    $newRange.Text = "Notifications"
}
